# Clear the "No. of Sites/bldg ..." breakdown columns (AB:AK) and the
# DIFFERENCE column (AM) for the data rows 2-18, while leaving the
# PREVIOUS ACCOMPLISHMENT column (AL) untouched, per the latest status
# accomplishment update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2:AK18").ClearContents()
$ws.Range("AM2:AM18").ClearContents()
